$d = $word.ActiveDocument

# --- 1. "900 to 1200 word" paragraph: merge the three runs (removing the
#        gramStart/gramEnd proofErr markers around "900 to 1200 word") by
#        replacing the whole span with identical text in one Find/Replace
#        pass, which causes the run to be re-authored as a single run.
$r1 = $d.Content
$r1.Find.Execute(
    "quality and completeness of a 900 to 1200 word ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "quality and completeness of a 900 to 1200 word ", 2) | Out-Null

# --- 2. "screenshare i.e. loom.com" paragraph: merge the three runs
#        (removing the gramStart/gramEnd proofErr markers around "i.e.")
#        the same way.
$r2 = $d.Content
$r2.Find.Execute(
    ", or a voice over in the slide file, screenshare i.e. loom.com or shared in a similarly appropriate manner. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ", or a voice over in the slide file, screenshare i.e. loom.com or shared in a similarly appropriate manner. ", 2) | Out-Null

# --- 3. Schedule table: shift the WallStreetBets/GameStop case ("Mar 8"
#        session) forward to "Mar 10". Do this one FIRST (before renaming
#        the "Mar 4" row to "Mar 8") so the Find below unambiguously hits
#        the original "Mar 8" cell, not a freshly produced one. The cell is
#        split into two runs ("Mar " / "10") the way Word does when only
#        the trailing digit(s) of the cell are retyped.
$r3 = $d.Content
$r3.Find.Execute("Mar 8", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$last3 = $d.Range($r3.End - 1, $r3.End)
$last3.Font.Bold = 1
$last3.Text = "10"
$new3 = $d.Range($last3.Start, $last3.End)
$new3.Font.Bold = 0

# --- 4. Schedule table: shift the prior session date "Mar 4" forward to
#        "Mar 8", split the same way.
$r4 = $d.Content
$r4.Find.Execute("Mar 4", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$last4 = $d.Range($r4.End - 1, $r4.End)
$last4.Font.Bold = 1
$last4.Text = "8"
$new4 = $d.Range($last4.Start, $last4.End)
$new4.Font.Bold = 0
